$wb = $excel.ActiveWorkbook

# --- Update the localization status text: "Ready for handoff" -> "In Translation" ---
# This shared string is used on the "Overview" sheet (columns E/F, the zh-cn/de-de
# status cells) as well as on the "zh-cn" and "de-de" detail sheets (column C,
# "Status"). Find/replace across every worksheet so all occurrences move together.
foreach ($sheet in $wb.Worksheets) {
    $found = $sheet.Cells.Find("Ready for handoff")
    if ($found -ne $null) {
        $firstAddress = $found.Address()
        do {
            $found.Value = "In Translation"
            $found = $sheet.Cells.FindNext($found)
        } while ($found -ne $null -and $found.Address() -ne $firstAddress)
    }
}

# --- Re-fit the status columns now that the text is shorter ---
# "In Translation" is narrower than "Ready for handoff", so the columns that hold
# the status values shrink accordingly when the report is regenerated.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
